$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Range("B2").Value = -0.2951083986270988
$ws1.Range("B3").Value = 16.37287915320643
$ws1.Range("B4").Value = 15.96882674611816
$ws1.Range("B5").Value = 15.75092415266201
$ws1.Range("B6").Value = 15.38674919193659
$ws1.Range("B7").Value = 16.10869972702038
$ws1.Range("B8").Value = 15.18691985173505
$ws1.Range("B9").Value = 14.96656560824961
$ws1.Range("B10").Value = 14.76370525054857
$ws1.Range("B11").Value = 14.15058926484286
$ws1.Range("B12").Value = 13.55394675651172
$ws1.Range("B13").Value = 13.04499173792274
$ws1.Range("B14").Value = 12.15534137340026
$ws1.Range("B15").Value = 11.85011782321278
$ws1.Range("B16").Value = 11.31159652291865
$ws1.Range("B17").Value = 10.68505103650835
$ws1.Range("B18").Value = 10.22169001469732
$ws1.Range("B19").Value = 9.900378632784362
$ws1.Range("B20").Value = 9.350810685155519
$ws1.Range("B21").Value = 9.211369153863844
$ws1.Range("B22").Value = 8.811704178438253
$ws1.Range("B23").Value = 9.202377947627543
$ws1.Range("B24").Value = 9.043505529180418
$ws1.Range("B25").Value = 8.99977847828092
$ws1.Range("B26").Value = 9.210595235470345
$ws1.Range("B27").Value = 9.030898630616646
$ws1.Range("B28").Value = 8.918249104323483

$ws2 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws2.Range("B2").Value = -0.2548739821090755
$ws2.Range("B3").Value = 16.41311356972444
$ws2.Range("B4").Value = 16.00906116263618
$ws2.Range("B5").Value = 15.79115856918003
$ws2.Range("B6").Value = 15.42698360845461
$ws2.Range("B7").Value = 16.1489341435384
$ws2.Range("B8").Value = 15.22715426825308
$ws2.Range("B9").Value = 15.00680002476764
$ws2.Range("B10").Value = 14.8039396670666
$ws2.Range("B11").Value = 14.19082368136089
$ws2.Range("B12").Value = 13.59418117302975
$ws2.Range("B13").Value = 13.08522615444076
$ws2.Range("B14").Value = 12.19557578991829
$ws2.Range("B15").Value = 11.89035223973081
$ws2.Range("B16").Value = 11.35183093943668
$ws2.Range("B17").Value = 10.72528545302638
$ws2.Range("B18").Value = 10.26192443121534
$ws2.Range("B19").Value = 9.940613049302385
$ws2.Range("B20").Value = 9.391045101673543
$ws2.Range("B21").Value = 9.251603570381867
$ws2.Range("B22").Value = 8.851938594956277
$ws2.Range("B23").Value = 9.242612364145566
$ws2.Range("B24").Value = 9.083739945698442
$ws2.Range("B25").Value = 9.040012894798943
$ws2.Range("B26").Value = 9.210595235470345
$ws2.Range("B27").Value = 9.071133047134669
$ws2.Range("B28").Value = 8.958483520841506
$ws2.Range("B29").Value = 8.740295205419871
$ws2.Range("B30").Value = 8.692727588217794
$ws2.Range("B31").Value = 8.592107686228566
$ws2.Range("B32").Value = 8.624854193911407
$ws2.Range("B33").Value = 8.442980026087612
